# Refresh cryptos list (Coin/Link/Price/Volume) with latest scraped values.
# Mirrors the upstream GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.857.18"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "1.635.02"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  +0.86%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.33"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.67"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("D12").Value = "1.870.11"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").Value = "1.640.49"
$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("E14").Value = "  +3.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.41"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +6.48%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.85"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -1.53%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.873.35"
$ws.Range("E17").Value = "  -0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.57"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.40"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("D20").Value = "0.0₃0702"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("E21").Value = "  +0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.88"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +2.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.71"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.52"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.109"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.62"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("E29").Value = "  +0.79%  "

$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.96%  "

$ws.Range("D34").Value = "1.425.65"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("E35").Value = "  +2.88%  "

$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("E37").Value = "  -3.10%  "

$ws.Range("E38").Value = "  +1.51%  "

$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.20"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +10.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.560"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.832"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  -0.79%  "

$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("E45").Value = "  +0.85%  "

$ws.Range("E46").Value = "  -1.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.40"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.59%  "

$ws.Range("D48").Value = "1.777.16"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.93"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -8.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.75"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +4.84%  "

$ws.Range("D51").Value = "0.0₆0109"
$ws.Range("E51").Value = "  +2.72%  "
